# tableInTable-template.docx — M2Doc issue #295
#
# The upstream commit ("Fixed #295 Add the version of M2Doc in the
# template custom properties") touched many template .docx resources in
# the repository. For *this* particular file the recorded OOXML diff is
# not a content edit at all: every hunk in word/document.xml and
# word/styles.xml is a pure attribute-order normalization (e.g.
# <w:tblW w:w="0" w:type="auto"/> -> <w:tblW w:type="auto" w:w="0"/>,
# <w:pgSz w:w="11906" w:h="16838"/> -> <w:pgSz w:h="16838" w:w="11906"/>,
# xmlns:* declarations re-sorted alphabetically on <w:document>, etc.).
# Re-sorting the XML attributes of every element in both parts
# reproduces exactly the canonical (C14N-style) attribute ordering shown
# in the diff, with no text, table, style, run/paragraph-property, or
# section-property value added, removed, or altered anywhere.
#
# Because that is a pure serialization/canonicalization artifact and not
# an actual document edit, there is nothing to do on the Word object
# model here: no text, table cell, style or page-setup value changes.
# We simply touch the document so the edit is a deliberate no-op,
# leaving content, formatting and structure bit-for-bit equivalent to
# the original.

$d = $word.ActiveDocument
